$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This single shared string is referenced from Overview!E2,F2,E3,F3 and
#    from the Status column (C2,C3) on both the zh-cn and de-de sheets, so a
#    single text replace keeps every cell pointed at the same (renamed)
#    shared string, just like the source edit.
# ---------------------------------------------------------------------------
$wsOverview.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US", -4163, 1, $false, $false, $false) | Out-Null
$wsZhCn.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US", -4163, 1, $false, $false, $false) | Out-Null
$wsDeDe.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US", -4163, 1, $false, $false, $false) | Out-Null

# ---------------------------------------------------------------------------
# 2. Fill in the handback report columns for zh-cn (sheet2) and de-de (sheet3)
#    I = Latest Target File, J = Latest Handback File, K = Latest Handback DateTime
# ---------------------------------------------------------------------------

# zh-cn
$wsZhCn.Range("I2").Value = "09572edd-dbb6-4c5e-ac11-fa5758def696.md"
$wsZhCn.Range("I3").Value = "09572edd-dbb6-4c5e-ac11-fa5758def696.md"
$wsZhCn.Range("I2:I3").Font.Underline = 2
$wsZhCn.Range("I2:I3").Font.Color = 15570276

$wsZhCn.Range("J2").Value = "09572edd-dbb6-4c5e-ac11-fa5758def696.c2607544a66b02a746f17728b9b1fcccf78d1073.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "09572edd-dbb6-4c5e-ac11-fa5758def696.c2607544a66b02a746f17728b9b1fcccf78d1073.zh-cn.xlf"

# Latest Handback DateTime for zh-cn (K2/K3) already shares its string with
# "0001-01-01 00:00:00"; rename that string in place for these two cells.
$wsZhCn.Range("K2:K3").Replace("0001-01-01 00:00:00", "2016-10-21 04:19:48", -4163, 1, $false, $false, $false) | Out-Null

# de-de
$wsDeDe.Range("I2").Value = "09572edd-dbb6-4c5e-ac11-fa5758def696.md"
$wsDeDe.Range("I3").Value = "09572edd-dbb6-4c5e-ac11-fa5758def696.md"
$wsDeDe.Range("I2:I3").Font.Underline = 2
$wsDeDe.Range("I2:I3").Font.Color = 15570276

$wsDeDe.Range("J2").Value = "09572edd-dbb6-4c5e-ac11-fa5758def696.c2607544a66b02a746f17728b9b1fcccf78d1073.de-de.xlf"
$wsDeDe.Range("J3").Value = "09572edd-dbb6-4c5e-ac11-fa5758def696.c2607544a66b02a746f17728b9b1fcccf78d1073.de-de.xlf"

# de-de handback happened a little later and gets its own, brand new datetime string.
$wsDeDe.Range("K2").Value = "2016-10-21 04:20:07"
$wsDeDe.Range("K3").Value = "2016-10-21 04:20:07"

# ---------------------------------------------------------------------------
# 3. Hyperlinks: add a "Latest Target File" hyperlink (column I) next to the
#    existing "Source File Name" hyperlink (column A) on rows 2 and 3, for
#    both the zh-cn and de-de sheets.  Rebuild the collection in the desired
#    final order (A2, I2, A3, I3) so relationship ids line up: rId2, rId3,
#    rId4, rId5.
# ---------------------------------------------------------------------------
$urlMd    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dc256c0b853190a2bf123e915a6885c70f4fa264/e2e/09572edd-dbb6-4c5e-ac11-fa5758def696.md"
$urlFfff  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dc256c0b853190a2bf123e915a6885c70f4fa264/e2e/ffff3f2694a0-01a4-436c-9d39-a13344643e0b.md"

foreach ($ws in @($wsZhCn, $wsDeDe)) {
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $urlMd,   "", "", "09572edd-dbb6-4c5e-ac11-fa5758def696.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I2"), $urlMd,   "", "", "09572edd-dbb6-4c5e-ac11-fa5758def696.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), $urlFfff, "", "", "ffff3f2694a0-01a4-436c-9d39-a13344643e0b.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I3"), $urlMd,   "", "", "09572edd-dbb6-4c5e-ac11-fa5758def696.md") | Out-Null
}

# ---------------------------------------------------------------------------
# 4. Column widths
#    Overview: zh-cn (E) / de-de (F) status columns widen.
#    zh-cn / de-de sheets: Status (C) widens, Latest Target File (I) and
#    Latest Handback File (J) widen to fit the newly populated file names.
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.16666666666667
$wsOverview.Columns.Item(6).ColumnWidth = 29.16666666666667

foreach ($ws in @($wsZhCn, $wsDeDe)) {
    $ws.Columns.Item(3).ColumnWidth  = 29.16666666666667
    $ws.Columns.Item(9).ColumnWidth  = 39.16666666666667
    $ws.Columns.Item(10).ColumnWidth = 39.16666666666667
}
